$d = $word.ActiveDocument

# 1) Title: "Finance-Loan Application Project" - shrink from 32 -> 28 half-points (16pt -> 14pt)
$pTitle = $d.Paragraphs(1)
$pTitle.Range.Font.Size = 14
$pTitle.Range.Font.SizeBi = 14

# 2) Image paragraph (centered drawing) - add sz=28/szCs=28 to paragraph mark + run
$pImage = $d.Paragraphs(7)
$pImage.Range.Font.Size = 14
$pImage.Range.Font.SizeBi = 14

# 3) "Testing data:" heading - add sz=28/szCs=28
$pHeading = $d.Paragraphs(8)
$pHeading.Range.Font.Size = 14
$pHeading.Range.Font.SizeBi = 14

# 4) JSON sample block (paragraphs 9 "{" .. 25 "}") - add sz=28/szCs=28 to every paragraph
for ($i = 9; $i -le 25; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
}

# 5) Replace the two trailing blank paragraphs with the new service-port lines.
#    Paragraph 26 already exists (blank) - fill it in directly and format it.
$p26 = $d.Paragraphs(26)
$p26.Range.InsertAfter("API Gateway:8080")
$p26 = $d.Paragraphs(26)
$p26.Range.Font.Size = 14
$p26.Range.Font.SizeBi = 14

# Split a fresh paragraph off paragraph 26 for each subsequent line - the newly minted
# paragraph mark inherits formatting from the range it was split from.
$p26.Range.InsertParagraphAfter()
$p27 = $d.Paragraphs(27)
$p27.Range.InsertAfter("Customer Port:3001")
$p27 = $d.Paragraphs(27)
$p27.Range.Font.Size = 14
$p27.Range.Font.SizeBi = 14

$p27.Range.InsertParagraphAfter()
$p28 = $d.Paragraphs(28)
$p28.Range.InsertAfter("Admin port :3002")
$p28 = $d.Paragraphs(28)
$p28.Range.Font.Size = 14
$p28.Range.Font.SizeBi = 14

$p28.Range.InsertParagraphAfter()
$p29 = $d.Paragraphs(29)
$p29.Range.InsertAfter("LoanService:3003")
$p29 = $d.Paragraphs(29)
$p29.Range.Font.Size = 14
$p29.Range.Font.SizeBi = 14

# 6) Final trailing blank paragraph (originally the very last, wholly-empty paragraph in
#    the document) needs sz=28/szCs=28 on its paragraph mark too. A completely empty
#    paragraph's Range silently refuses direct Font writes, so briefly give it a
#    placeholder character to anchor the formatting, apply the size, then remove the
#    placeholder again - leaving just the formatted, empty paragraph mark behind.
$p30 = $d.Paragraphs(30)
$p30.Range.InsertAfter("x")
$p30 = $d.Paragraphs(30)
$p30.Range.Font.Size = 14
$p30.Range.Font.SizeBi = 14
$p30 = $d.Paragraphs(30)
$clearRng = $d.Range($p30.Range.Start, $p30.Range.Start + 1)
$clearRng.Text = ""

Write-Output "Pagination + LoanService update applied."
